$wb = $excel.ActiveWorkbook

# --- TestStep sheet: E1 content change ("dummy" -> "Release") ---
$wsTestStep = $wb.Worksheets.Item("TestStep")
$wsTestStep.Range("E1").Value = "Release"

# --- TestCaseSequence sheet: add new "LineNumbers" column (H) ---
$wsSeq = $wb.Worksheets.Item("TestCaseSequence")
$wsSeq.Range("H1").Value = "LineNumbers"
$wsSeq.Range("H2").Value = "1:5, 6, 9, 11"

# Column widths to roughly match the authored layout
$wsSeq.Columns("C").ColumnWidth = 22.666666666666668
$wsSeq.Columns("H").ColumnWidth = 11.8

# --- Selection / active-sheet bookkeeping to mirror the authored session ---
$wsExport = $wb.Worksheets.Item("ExportFieldList")
$null = $wsExport.Range("A1").Select()

$wsTestCase = $wb.Worksheets.Item("TestCase")
$null = $wsTestCase.Range("C2").Select()

$null = $wsTestStep.Range("E2").Select()

$wsExec = $wb.Worksheets.Item("TestStepExecution")
$null = $wsExec.Range("E8").Select()

$null = $wsSeq.Range("H2").Select()
$null = $wsSeq.Activate()
